# "added new round data and Tee Fairway column"
#
# 1. Rounds sheet gets a new entry for the 06/02/2025 round.
# 2. The four existing per-round scorecard sheets get a new "Tee Fairway"
#    column inserted at C (pushing Fairway Hits/Chips/Putts to D/E/F); the
#    three later sheets' headers also pick up a slightly different font.
# 3. A brand new scorecard sheet "CGC-W-06022025" is added at the end with
#    full hole-by-hole data, including the new Tee Fairway column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add a new round entry to the "Rounds" sheet.
# ---------------------------------------------------------------------------
$roundsWs = $wb.Worksheets.Item("Rounds")
$roundsWs.Cells.Item(6, 1).Value = "CGC-W-06022025"
$roundsWs.Cells.Item(6, 2).Value = "CGC-W"

# Copy the date format from the row above so the new date cell reuses the
# existing "m/d/yyyy" style instead of creating a new one.
$roundsWs.Cells.Item(5, 3).Copy()
$roundsWs.Cells.Item(6, 3).PasteSpecial(-4122)  # xlPasteFormats
$roundsWs.Cells.Item(6, 3).Value = (Get-Date -Year 2025 -Month 6 -Day 2).Date
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Insert a "Tee Fairway" column (new column C) into the four existing
#    per-round scorecard sheets, shifting Fairway Hits/Chips/Putts right.
# ---------------------------------------------------------------------------
$roundSheetNames = @("CGC-W-05222025", "CGC-W-05232025", "CGC-W-05302025", "CGC-W-05312025")

foreach ($name in $roundSheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(3).Insert()
    $ws.Cells.Item(1, 3).Value = "Tee Fairway"
}

# The three later sheets (CGC-W-05232025, CGC-W-05302025, CGC-W-05312025)
# pick up a slightly different header font; CGC-W-05222025 keeps the default.
$styledSheetNames = @("CGC-W-05232025", "CGC-W-05302025", "CGC-W-05312025")
foreach ($name in $styledSheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $headerRange = $ws.Range("C1:F1")
    $headerRange.Font.Color = 0
    $headerRange.Font.Size = 12
}

# ---------------------------------------------------------------------------
# 3. Add the new scorecard sheet "CGC-W-06022025" with full hole-by-hole data
#    (including the new Tee Fairway column).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "CGC-W-06022025"

$newSheet.Cells.Item(1, 1).Value = "Hole"
$newSheet.Cells.Item(1, 2).Value = "Score"
$newSheet.Cells.Item(1, 3).Value = "Tee Fairway"
$newSheet.Cells.Item(1, 4).Value = "Fairway Hits"
$newSheet.Cells.Item(1, 5).Value = "Chips"
$newSheet.Cells.Item(1, 6).Value = "Putts"

$holeData = @(
    @(1, 9, "Yes", 0, 1, 2),
    @(2, 7, "Yes", 0, 0, 2),
    @(3, 4, "No",  0, 1, 2),
    @(4, 8, "Yes", 0, 4, 2),
    @(5, 6, "Yes", 1, 2, 2),
    @(6, 7, "No",  1, 1, 2),
    @(7, 4, "No",  0, 1, 2),
    @(8, 7, "No",  1, 1, 2),
    @(9, 8, "No",  2, 2, 2)
)

for ($i = 0; $i -lt $holeData.Length; $i++) {
    $r = $i + 2
    $row = $holeData[$i]
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
}

# ---------------------------------------------------------------------------
# 4. Selections / active sheet to match the after-state of the edit.
# ---------------------------------------------------------------------------
$roundsWs.Range("C7").Select() | Out-Null
$wb.Worksheets.Item("CGC-W-05222025").Range("F13").Select() | Out-Null
$wb.Worksheets.Item("CGC-W-05232025").Range("F2:F10").Select() | Out-Null
$wb.Worksheets.Item("CGC-W-05302025").Range("F2:F10").Select() | Out-Null
$wb.Worksheets.Item("CGC-W-05312025").Range("F2:F10").Select() | Out-Null

$newSheet.Activate() | Out-Null
$newSheet.Range("D11").Select() | Out-Null
